$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New header cells F1:J1, matching style of existing header cells ---
$ws.Range("B1").Copy()
$ws.Range("F1:J1").PasteSpecial(-4122)
$excel.CutCopyMode = 0

$ws.Range("F1").Value = "Sensitivity"
$ws.Range("G1").Value = "Specificity"
$ws.Range("H1").Value = "Precision"
$ws.Range("I1").Value = "Recall"
$ws.Range("J1").Value = "F1 Score"

# --- Update existing AUC / CI bound columns (C:E) with new values ---
$ws.Range("C2").Value = 0.8263805067323644
$ws.Range("D2").Value = 0.7883910013903975
$ws.Range("E2").Value = 0.8630535951604434

$ws.Range("C3").Value = 0.8042571250626939
$ws.Range("D3").Value = 0.7653966682008022
$ws.Range("E3").Value = 0.8417185094358814

$ws.Range("C4").Value = 0.7340867628984165
$ws.Range("D4").Value = 0.6961271644701137
$ws.Range("E4").Value = 0.7706563481787226

# --- New data columns F:J ---
$ws.Range("F2").Value = 0.9370249728555917
$ws.Range("G2").Value = 0.7157360406091371
$ws.Range("H2").Value = 0.9390642002176278
$ws.Range("I2").Value = 0.9370249728555917
$ws.Range("J2").Value = 0.9380434782608695

$ws.Range("F3").Value = 0.9435396308360477
$ws.Range("G3").Value = 0.6649746192893401
$ws.Range("H3").Value = 0.9294117647058824
$ws.Range("I3").Value = 0.9435396308360477
$ws.Range("J3").Value = 0.9364224137931034

$ws.Range("F4").Value = 0.9554831704668838
$ws.Range("G4").Value = 0.5126903553299492
$ws.Range("H4").Value = 0.9016393442622951
$ws.Range("I4").Value = 0.9554831704668838
$ws.Range("J4").Value = 0.9277807063784924
